# Correcting formula for check #4
# Recomputed structural-check outputs (selected member size, resistances, utilisations)
# for the strut rows in the processed results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 2915
$ws.Range("Y2").Value = "'286"
$ws.Range("Y2").Style = "Normal"
$ws.Range("Z2").Value = "'3202"
$ws.Range("Z2").Style = "Normal"
$ws.Range("AA2").Value = "457 x 191 x 74"
$ws.Range("AC2").Value = 6717
$ws.Range("AD2").Value = 5502
$ws.Range("AE2").Value = 1.1
$ws.Range("AG2").Value = 0.59
$ws.Range("AH2").Value = 3981
$ws.Range("AJ2").Value = 6717
$ws.Range("AK2").Value = 473000
$ws.Range("AL2").Value = 90639
$ws.Range("AP2").Value = 6608
$ws.Range("AR2").Value = 3358
$ws.Range("AS2").Value = 54082
$ws.Range("AT2").Value = 0.25
$ws.Range("AV2").Value = 0.99
$ws.Range("AW2").Value = 3322
$ws.Range("AX2").Value = 1.18
$ws.Range("AZ2").Value = 78153
$ws.Range("BD2").Value = 1085
$ws.Range("BE2").Value = 3310

# Row 3
$ws.Range("Q3").Value = 1546
$ws.Range("Y3").Value = "'154"
$ws.Range("Y3").Style = "Normal"
$ws.Range("Z3").Value = "'1700"
$ws.Range("Z3").Style = "Normal"
$ws.Range("AA3").Value = "305 x 165 x 40"
$ws.Range("AC3").Value = 3642
$ws.Range("AD3").Value = 5618
$ws.Range("AE3").Value = 0.8100000000000001
$ws.Range("AG3").Value = 0.79
$ws.Range("AH3").Value = 2887
$ws.Range("AJ3").Value = 3642
$ws.Range("AK3").Value = 256500
$ws.Range("AL3").Value = 169524
$ws.Range("AP3").Value = 3685
$ws.Range("AR3").Value = 1821
$ws.Range("AS3").Value = 24742
$ws.Range("AT3").Value = 0.27
$ws.Range("AV3").Value = 0.98
$ws.Range("AW3").Value = 1792
$ws.Range("AX3").Value = 0.85
$ws.Range("AZ3").Value = 169524
$ws.Range("BA3").Value = 52
$ws.Range("BB3").Value = 218
$ws.Range("BC3").Value = 206
$ws.Range("BD3").Value = 259
$ws.Range("BE3").Value = 1726

# Row 4
$ws.Range("Q4").Value = 4810
$ws.Range("Z4").Value = "'4810"
$ws.Range("Z4").Style = "Normal"
$ws.Range("AA4").Value = "610 x 229 x 113"
$ws.Range("AB4").Value = 345
$ws.Range("AC4").Value = 9936
$ws.Range("AD4").Value = 14424
$ws.Range("AE4").Value = 0.83
$ws.Range("AG4").Value = 0.78
$ws.Range("AH4").Value = 7732
$ws.Range("AI4").Value = 345
$ws.Range("AJ4").Value = 9936
$ws.Range("AK4").Value = 720000
$ws.Range("AL4").Value = 137970
$ws.Range("AP4").Value = 9785
$ws.Range("AQ4").Value = 345
$ws.Range("AR4").Value = 4968
$ws.Range("AS4").Value = 111079
$ws.Range("AT4").Value = 0.21
$ws.Range("AV4").Value = 1
$ws.Range("AW4").Value = 4955
$ws.Range("AX4").Value = 1.56
$ws.Range("AZ4").Value = 118964
$ws.Range("BD4").Value = 1182
$ws.Range("BE4").Value = 4928

# Row 5
$ws.Range("Q5").Value = 2551
$ws.Range("Z5").Value = "'2551"
$ws.Range("Z5").Style = "Normal"
$ws.Range("AA5").Value = "457 x 152 x 60"
$ws.Range("AC5").Value = 5410
$ws.Range("AD5").Value = 16853
$ws.Range("AE5").Value = 0.57
$ws.Range("AG5").Value = 0.9
$ws.Range("AH5").Value = 4881
$ws.Range("AJ5").Value = 5410
$ws.Range("AK5").Value = 381000
$ws.Range("AL5").Value = 251807
$ws.Range("AP5").Value = 5473
$ws.Range("AR5").Value = 2705
$ws.Range("AS5").Value = 25746
$ws.Range("AT5").Value = 0.32
$ws.Range("AW5").Value = 2629
$ws.Range("AX5").Value = 1.04
$ws.Range("AZ5").Value = 251807
$ws.Range("BA5").Value = 52
$ws.Range("BB5").Value = 218
$ws.Range("BC5").Value = 206
$ws.Range("BD5").Value = 283
$ws.Range("BE5").Value = 2579

# Row 6
$ws.Range("Q6").Value = 6885
$ws.Range("Z6").Value = "'6885"
$ws.Range("Z6").Style = "Normal"
$ws.Range("AA6").Value = "457 x 191 x 161"
$ws.Range("AB6").Value = 345
$ws.Range("AC6").Value = 14214
$ws.Range("AD6").Value = 13185
$ws.Range("AE6").Value = 1.04
$ws.Range("AG6").Value = 0.64
$ws.Range("AH6").Value = 9080
$ws.Range("AI6").Value = 345
$ws.Range("AJ6").Value = 14214
$ws.Range("AK6").Value = 1030000
$ws.Range("AL6").Value = 197374
$ws.Range("AP6").Value = 13998
$ws.Range("AQ6").Value = 345
$ws.Range("AR6").Value = 7107
$ws.Range("AS6").Value = 137635
$ws.Range("AT6").Value = 0.23
$ws.Range("AV6").Value = 0.99
$ws.Range("AW6").Value = 7064
$ws.Range("AX6").Value = 2.03
$ws.Range("AZ6").Value = 170185
$ws.Range("BA6").Value = 53
$ws.Range("BB6").Value = 445
$ws.Range("BC6").Value = 833
$ws.Range("BD6").Value = 1326
$ws.Range("BE6").Value = 7018

# Row 7
$ws.Range("Q7").Value = 3651
$ws.Range("Z7").Value = "'3651"
$ws.Range("Z7").Style = "Normal"
$ws.Range("AA7").Value = "533 x 210 x 82"
$ws.Range("AC7").Value = 7455
$ws.Range("AD7").Value = 31393
$ws.Range("AE7").Value = 0.49
$ws.Range("AG7").Value = 0.93
$ws.Range("AH7").Value = 6919
$ws.Range("AJ7").Value = 7455
$ws.Range("AK7").Value = 525000
$ws.Range("AL7").Value = 346978
$ws.Range("AP7").Value = 7541
$ws.Range("AR7").Value = 3728
$ws.Range("AS7").Value = 65093
$ws.Range("AT7").Value = 0.24
$ws.Range("AV7").Value = 0.99
$ws.Range("AW7").Value = 3695
$ws.Range("AX7").Value = 1.26
$ws.Range("AZ7").Value = 346978
$ws.Range("BD7").Value = 315
$ws.Range("BE7").Value = 3682

# Row 8
$ws.Range("Q8").Value = 5704
$ws.Range("Z8").Value = "'5704"
$ws.Range("Z8").Style = "Normal"
$ws.Range("AA8").Value = "762 x 267 x 134"
$ws.Range("AC8").Value = 12141
$ws.Range("AD8").Value = 24949
$ws.Range("AE8").Value = 0.7
$ws.Range("AG8").Value = 0.85
$ws.Range("AH8").Value = 10306
$ws.Range("AJ8").Value = 12141
$ws.Range("AK8").Value = 855000
$ws.Range("AL8").Value = 163839
$ws.Range("AP8").Value = 11946
$ws.Range("AR8").Value = 6070
$ws.Range("AS8").Value = 155122
$ws.Range("AT8").Value = 0.2
$ws.Range("AV8").Value = 1
$ws.Range("AW8").Value = 6073
$ws.Range("AX8").Value = 1.77
$ws.Range("AZ8").Value = 141270
$ws.Range("BD8").Value = 1236
$ws.Range("BE8").Value = 5828

# Row 9
$ws.Range("Q9").Value = 3025
$ws.Range("Z9").Value = "'3025"
$ws.Range("Z9").Style = "Normal"
$ws.Range("AA9").Value = "457 x 191 x 74"
$ws.Range("AC9").Value = 6717
$ws.Range("AD9").Value = 22008
$ws.Range("AE9").Value = 0.55
$ws.Range("AG9").Value = 0.91
$ws.Range("AH9").Value = 6093
$ws.Range("AJ9").Value = 6717
$ws.Range("AK9").Value = 473000
$ws.Range("AL9").Value = 312611
$ws.Range("AP9").Value = 6795
$ws.Range("AR9").Value = 3358
$ws.Range("AS9").Value = 54082
$ws.Range("AT9").Value = 0.25
$ws.Range("AV9").Value = 0.99
$ws.Range("AW9").Value = 3322
$ws.Range("AX9").Value = 1.18
$ws.Range("AZ9").Value = 312611
$ws.Range("BD9").Value = 296
$ws.Range("BE9").Value = 3055

# Row 10
$ws.Range("Q10").Value = 3356
$ws.Range("Z10").Value = "'3356"
$ws.Range("Z10").Style = "Normal"
$ws.Range("AA10").Value = "610 x 178 x 82"
$ws.Range("AC10").Value = 7384
$ws.Range("AD10").Value = 9236
$ws.Range("AE10").Value = 0.89
$ws.Range("AG10").Value = 0.74
$ws.Range("AH10").Value = 5448
$ws.Range("AJ10").Value = 7384
$ws.Range("AK10").Value = 520000
$ws.Range("AL10").Value = 99645
$ws.Range("AP10").Value = 7265
$ws.Range("AR10").Value = 3692
$ws.Range("AS10").Value = 39185
$ws.Range("AT10").Value = 0.31
$ws.Range("AV10").Value = 0.98
$ws.Range("AW10").Value = 3603
$ws.Range("AX10").Value = 1.26
$ws.Range("AZ10").Value = 85918
$ws.Range("BD10").Value = 1090
$ws.Range("BE10").Value = 3465

# Row 11
$ws.Range("Q11").Value = 1780
$ws.Range("Z11").Value = "'1780"
$ws.Range("Z11").Style = "Normal"
$ws.Range("AA11").Value = "305 x 127 x 42"
$ws.Range("AC11").Value = 3791
$ws.Range("AD11").Value = 5419
$ws.Range("AE11").Value = 0.84
$ws.Range("AG11").Value = 0.77
$ws.Range("AH11").Value = 2935
$ws.Range("AJ11").Value = 3791
$ws.Range("AK11").Value = 267000
$ws.Range("AL11").Value = 176463
$ws.Range("AP11").Value = 3835
$ws.Range("AR11").Value = 1896
$ws.Range("AS11").Value = 12598
$ws.Range("AT11").Value = 0.39
$ws.Range("AV11").Value = 0.96
$ws.Range("AW11").Value = 1812
$ws.Range("AX11").Value = 0.87
$ws.Range("AZ11").Value = 176463
$ws.Range("BA11").Value = 52
$ws.Range("BB11").Value = 218
$ws.Range("BC11").Value = 206
$ws.Range("BD11").Value = 261
$ws.Range("BE11").Value = 1806
